# Updated cryptos list on Fri Mar 24 22:59:24 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for rows 2-51 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.399.30"
$ws.Range("E2").Value = "  -3.16%  "
$ws.Range("D3").Value = "1.749.29"
$ws.Range("E3").Value = "  -3.47%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'322.51"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4242"
$ws.Range("E7").Value = "  -4.93%  "
$ws.Range("D8").Value = "'0.3600"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "'0.07491"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").Value = "'42.12"
$ws.Range("E10").Value = "  -6.28%  "
$ws.Range("D11").Value = "'1.097"
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -6.45%  "
$ws.Range("D14").Value = "'6.026"
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("D15").Value = "'7.218"
$ws.Range("E15").Value = "  -4.73%  "
$ws.Range("D16").Value = "1.752.24"
$ws.Range("E16").Value = "  -5.33%  "
$ws.Range("D17").Value = "'92.72"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "'0.00001067"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").Value = "'0.06378"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'17.05"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").Value = "'5.893"
$ws.Range("E22").Value = "  -5.37%  "
$ws.Range("D23").Value = "27.445.50"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").Value = "'11.21"
$ws.Range("E24").Value = "  -4.05%  "
$ws.Range("D25").Value = "'2.092"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").Value = "'161.55"
$ws.Range("E26").Value = "  +3.63%  "
$ws.Range("D27").Value = "'20.25"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "1.949.83"
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("D29").Value = "'2.130"
$ws.Range("E29").Value = "  -7.74%  "
$ws.Range("D30").Value = "'123.73"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "'1.100"
$ws.Range("E31").Value = "  -8.25%  "
$ws.Range("D32").Value = "'3.647"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "'5.527"
$ws.Range("E33").Value = "  -6.42%  "
$ws.Range("D34").Value = "'0.08871"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").Value = "'12.20"
$ws.Range("E35").Value = "  -6.41%  "
$ws.Range("D36").Value = "'0.02284"
$ws.Range("E36").Value = "  -3.02%  "
$ws.Range("D37").Value = "'0.2097"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").Value = "'0.05999"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("D39").Value = "'0.6330"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").Value = "'4.935"
$ws.Range("E40").Value = "  -4.70%  "
$ws.Range("D41").Value = "'1.184"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "'7.887"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").Value = "'1.386"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").Value = "'13.40"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("D46").Value = "'0.5872"
$ws.Range("E46").Value = "  -3.55%  "
$ws.Range("D47").Value = "'3.693"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("D48").Value = "'1.970"
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").Value = "'122.34"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").Value = "'0.06808"
$ws.Range("E51").Value = "  -2.53%  "
